$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.203.54'
$ws.Range("E2").Value = '  -1.12%  '
$ws.Range("D3").Value = '1.659.36'
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5159'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.38%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2640'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06272'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07758'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.48%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.478'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.08%  '
$ws.Range("D13").Value = '1.658.81'
$ws.Range("E13").Value = '  -0.79%  '
$ws.Range("D14").Value = '1.886.49'
$ws.Range("E14").Value = '  -0.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5457'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.94%  '
$ws.Range("D16").Value = '0.0₅8128'
$ws.Range("E16").Value = '  -2.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.83'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.02%  '
$ws.Range("D18").Value = '26.213.98'
$ws.Range("E18").Value = '  -1.10%  '
$ws.Range("E20").Value = '  -3.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '192.17'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.51%  '
$ws.Range("E22").Value = '  -2.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.992'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '139.64'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1221'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.281'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.30%  '
$ws.Range("E28").Value = '  -1.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.440'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05935'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.274'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.545'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.271'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.582'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9613'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.428'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.770'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5677'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.48%  '
$ws.Range("E39").Value = '  +0.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01592'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8564'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.10%  '
$ws.Range("D43").Value = '1.011.79'
$ws.Range("E43").Value = '  -7.32%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.64'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("D45").Value = '1.801.41'
$ws.Range("E45").Value = '  -0.81%  '
$ws.Range("D46").Value = '0.0₈110'
$ws.Range("E46").Value = '  +2.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.47'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.002'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.019'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.45%  '
$ws.Range("E50").Value = '  -0.55%  '
$ws.Range("E51").Value = '  -0.74%  '
